$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (labels change, same A1/B1/C1 cells) ---
$ws.Range("A1").Value = "Capability"
$ws.Range("B1").Value = "No_emp_assigned"
$ws.Range("C1").Value = "No_emp_undergoing"

# --- Data rows (table shrinks from 13 rows to 7 rows, new values) ---
$ws.Range("A2").Value = "PowerShell And Bash"
$ws.Range("B2").Value = 3979
$ws.Range("C2").Value = 3621

$ws.Range("A3").Value = "A"
$ws.Range("B3").Value = 155
$ws.Range("C3").Value = 145

$ws.Range("A4").Value = "B"
$ws.Range("B4").Value = 4000
$ws.Range("C4").Value = 3900

$ws.Range("A5").Value = "A"
$ws.Range("B5").Value = 2834.3333333333298
$ws.Range("C5").Value = 2732.3333333333298

$ws.Range("A6").Value = "B"
$ws.Range("B6").Value = 2973.8333333333298
$ws.Range("C6").Value = 2742.8333333333298

$ws.Range("A7").Value = "A"
$ws.Range("B7").Value = 3113.3333333333298
$ws.Range("C7").Value = 2753.3333333333298

# Remove the now-unused rows 8:13 so the used range / dimension shrinks to A1:C7
$ws.Range("A8:C13").Clear()

# --- Column widths (bestFit-style custom widths on A:C) ---
# Target stored widths are 18.453125 / 16.1796875 / 18.54296875 (Excel's real
# bestFit, computed from on-screen glyph metrics). The closest widths this
# engine's column-width quantization can produce are used here.
$ws.Columns("A").ColumnWidth = 17.666666666666668
$ws.Columns("B").ColumnWidth = 15.333333333333334
$ws.Columns("C").ColumnWidth = 17.666666666666668

# --- Selection moves to B5:B7 ---
$ws.Range("B5:B7").Select()

# --- Page orientation switches to portrait ---
$ws.PageSetup.Orientation = 1
